# Apply the edit described by the diff:
#  1. Insert a new row before row 16 (end of the "tclif" block) holding a new
#     "fb"-style result row (Network name left blank, like the other repeated
#     rows in that block) with Final Accuracy = 92.875.
#  2. Insert a new "fb" row right after the "resnet" row of the "lif" block
#     (originally row 19, now row 20 after step 1) with Final Accuracy = 88.5.
#  3. Insert a new "fb" row right after the "alexnet" row of the "alif" block
#     (originally row 22, now row 24 after steps 1-2) with Final Accuracy = 22.9375.
#
# Every insert shifts the following rows down by one, so later row numbers
# account for the prior inserts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: new row at 16, appended to the "tclif" block ---
$ws.Rows("16:16").Insert()
$ws.Range("A16").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "Valve"
$ws.Range("D16").Value = "adam"
$ws.Range("E16").Value = 0.0005
$ws.Range("F16").Value = 256
$ws.Range("G16").Value = 200
$ws.Range("H16").Value = 92.875

# --- Step 2: new "fb" row at 21, appended to the "lif" block ---
$ws.Rows("21:21").Insert()
$ws.Range("A21").Value = ""
$ws.Range("B21").Value = "fb"
$ws.Range("C21").Value = "Valve"
$ws.Range("D21").Value = "adam"
$ws.Range("E21").Value = 0.0005
$ws.Range("F21").Value = 256
$ws.Range("G21").Value = 200
$ws.Range("H21").Value = 88.5

# --- Step 3: new "fb" row at 25, appended to the "alif" block ---
$ws.Rows("25:25").Insert()
$ws.Range("A25").Value = ""
$ws.Range("B25").Value = "fb"
$ws.Range("C25").Value = "Valve"
$ws.Range("D25").Value = "adam"
$ws.Range("E25").Value = 0.0005
$ws.Range("F25").Value = 256
$ws.Range("G25").Value = 200
$ws.Range("H25").Value = 22.9375
